$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new MCH122-1 record as row 2 (matches headers in row 1:
# identifier | alternativeIdentifiers | title | date_s | levelOfDescription | extentAndMedium | notes | file_path)
$ws.Range("A2").Value = "MCH122-1"
$ws.Range("C2").Value = "END CONSCRIPTION CAMPAIGN INCLUDING GRAPHICS"
$ws.Range("E2").Value = "Series"
$ws.Range("F2").Value = "1 Box"
$ws.Range("G2").Value = "LOCATION: 21D | GRAP COUNT NUMER: NONE"

# Give the new row the same "data row" look (non-bold Calibri 10pt, automatic
# text color) as distinct from the bold header row above it. Column B
# (alternativeIdentifiers) is left completely untouched/blank for this record,
# so it is formatted separately from the rest of the row.
$ws.Range("A2").Font.Name = "Calibri"
$ws.Range("A2").Font.ThemeColor = 1
$ws.Range("C2:H2").Font.Name = "Calibri"
$ws.Range("C2:H2").Font.ThemeColor = 1

# Re-establish the frozen header pane / selection on the new active row,
# matching the original frozen-top-row view.
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("A2:H2").Select()
